$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 9797.65
$ws.Range("B6").Value = 9920.67
$ws.Range("C6").Value = 19.36
$ws.Range("D6").Value = 19.12
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -1.24
$ws.Range("G6").Value = 42612.674803240741
$ws.Range("G6").NumberFormat = "m/d/yy h:mm"
$ws.Range("H6").Value = $false
